$d = $word.ActiveDocument

# The paragraph currently contains one big run:
#   ". Duis urna justo, vehicula vitae ultricies vel, congue at sem. Fusce
#    turpis turpis, ... Donec iaculis sed urna eget pharetra. "
# The target splits it into seven runs:
#   ". " / "Duis urna " / "justo"(bold) / ", " / "vehicula"(italic) /
#   " vitae ultricies vel, congue at sem." / " Fusce turpis turpis...pharetra. "
# None of these operations change the visible text, only formatting / run
# boundaries.

# 1. Split off ". " from "Duis urna " (both stay unformatted) by toggling
#    Bold on/off on "Duis urna " -- this creates a run boundary without
#    altering the rendered formatting.
$rDuisUrna = $d.Content.Duplicate
$rDuisUrna.Find.Execute("Duis urna justo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rDuisUrna.End = $rDuisUrna.Start + 10
$rDuisUrna.Bold = 1
$rDuisUrna.Bold = 0

# 2. Bold the word "justo" (search on the unique phrase "urna justo" and
#    narrow the range so only "justo" itself is affected).
$rJusto = $d.Content.Duplicate
$rJusto.Find.Execute("urna justo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rJusto.Start = $rJusto.Start + 5
$rJusto.Bold = 1

# 3. Italicize the word "vehicula".
$rVehicula = $d.Content.Duplicate
$rVehicula.Find.Execute("vehicula", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rVehicula.Italic = 1

# 4. Introduce a run break right before " Fusce turpis turpis, ..." without
#    changing any visible formatting there: toggle Bold on then back off on
#    that trailing range, which leaves the text/formatting unchanged but
#    keeps it as a distinct run from the preceding sentence.
$rRest = $d.Content.Duplicate
$rRest.Find.Execute(" Fusce turpis turpis, aliquet id pulvinar aliquam, iaculis non elit. Nulla feugiat lectus nulla, in dictum ipsum cursus ac. Quisque at odio neque. Sed ac tortor iaculis, bibendum leo ut, malesuada velit. Donec iaculis sed urna eget pharetra. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rRest.Bold = 1
$rRest.Bold = 0
